$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price figures stored as literal text (e.g. "1.000", "0.9999",
# "3.680", "25.991.02"). When such numeric-looking strings are assigned via
# .Value, Excel's automatic type detection silently converts them to Double
# (losing formatting / trailing zeros, e.g. "3.680" -> 3.68, "1.000" -> 1).
# To keep them as exact text -- matching the original inlineStr cells -- the
# Text number format ("@") is applied to each D cell immediately before its
# value is written.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.991.02"
$ws.Range("E2").Value = "  +0.12%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.742.16"
$ws.Range("E3").Value = "  +0.17%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.17%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.48"
$ws.Range("E5").Value = "  +3.75%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.09%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5051"
$ws.Range("E7").Value = "  -3.65%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2752"
$ws.Range("E8").Value = "  +0.73%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06189"
$ws.Range("E9").Value = "  +0.70%  "

# Row 10
$ws.Range("B10").Value = "WrappedEther"
$ws.Range("C10").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.742.59"
$ws.Range("E10").Value = "  +0.19%  "

# Row 11
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07246"
$ws.Range("E11").Value = "  +0.92%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.6543"
$ws.Range("E12").Value = "  +2.09%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.13"
$ws.Range("E13").Value = "  +0.25%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.686"
$ws.Range("E14").Value = "  +1.86%  "

# Row 15
$ws.Range("E15").Value = "  +0.20%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.0000"
$ws.Range("E16").Value = "  -0.12%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("E17").Value = "  -0.12%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.016.63"
$ws.Range("E18").Value = "  +0.10%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.94"
$ws.Range("E19").Value = "  +1.27%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006876"
$ws.Range("E20").Value = "  +2.00%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.969.63"
$ws.Range("E21").Value = "  +0.28%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.492"
$ws.Range("E22").Value = "  +3.55%  "

# Row 23
$ws.Range("E23").Value = "  +0.70%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.409"
$ws.Range("E24").Value = "  +3.02%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "135.82"
$ws.Range("E25").Value = "  -2.54%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.499"
$ws.Range("E26").Value = "  -1.01%  "

# Row 27
$ws.Range("E27").Value = "  +0.03%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.793"
$ws.Range("E28").Value = "  +1.44%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "105.45"
$ws.Range("E29").Value = "  -0.62%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.969"
$ws.Range("E30").Value = "  +1.97%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08173"
$ws.Range("E31").Value = "  -2.42%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.680"
$ws.Range("E32").Value = "  +1.45%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04693"
$ws.Range("E33").Value = "  +2.22%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.657"
$ws.Range("E34").Value = "  +0.14%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9955"
$ws.Range("E35").Value = "  +0.56%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6121"
$ws.Range("E36").Value = "  -1.94%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.754"
$ws.Range("E37").Value = "  +2.08%  "

# Row 38
$ws.Range("E38").Value = "  +1.30%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.932"
$ws.Range("E39").Value = "  -0.34%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9997"
$ws.Range("E40").Value = "  -0.09%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "101.04"
$ws.Range("E41").Value = "  +2.64%  "

# Row 42
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7950"
$ws.Range("E42").Value = "  +7.11%  "

# Row 43
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3921"
$ws.Range("E43").Value = "  +1.02%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.013"
$ws.Range("E44").Value = "  +1.32%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1162"
$ws.Range("E45").Value = "  +1.66%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.332"
$ws.Range("E46").Value = "  +1.45%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.65"
$ws.Range("E47").Value = "  +1.63%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05293"
$ws.Range("E48").Value = "  -0.19%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.78"
$ws.Range("E49").Value = "  +0.05%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.660"
$ws.Range("E50").Value = "  +1.57%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3480"
$ws.Range("E51").Value = "  +0.98%  "
